$wb = $excel.ActiveWorkbook

# Rename sheets: "rf" first (avoids name collision with existing "rf" sheet),
# then "arima" -> "rf", finally ensure second sheet is "lstm".
$wsArima = $wb.Sheets.Item(1)
$wsRf = $wb.Sheets.Item(2)
$wsRf.Name = "lstm"
$wsArima.Name = "rf"

# --- Sheet 1 ("rf", formerly "arima") updated metric values ---
$wsArima.Range("B2").Value = 0.5112219451371571
$wsArima.Range("C2").Value = 0.4314214463840399
$wsArima.Range("D2").Value = 0.5037406483790524
$wsArima.Range("E2").Value = 0.456359102244389
$wsArima.Range("F2").Value = 0.4488778054862843
$wsArima.Range("G2").Value = 0.5087281795511222
$wsArima.Range("H2").Value = 0.4987531172069826
$wsArima.Range("I2").Value = 0.5561097256857855
$wsArima.Range("B3").Value = 0.6756756756756757
$wsArima.Range("C3").Value = 0.4270833333333333
$wsArima.Range("D3").Value = 0.4758620689655172
$wsArima.Range("E3").Value = 0.3975155279503105
$wsArima.Range("F3").Value = 0.4113924050632912
$wsArima.Range("G3").Value = 0.4338235294117647
$wsArima.Range("H3").Value = 0.525
$wsArima.Range("I3").Value = 0.5571428571428572
$wsArima.Range("B4").Value = 78.40334261017935
$wsArima.Range("C4").Value = 0.3796945101752504
$wsArima.Range("D4").Value = 0.3615427946176608
$wsArima.Range("E4").Value = 0.008542353855724407
$wsArima.Range("F4").Value = 7.659918336857181
$wsArima.Range("G4").Value = 27.90679897716881
$wsArima.Range("H4").Value = 10741.31161794926
$wsArima.Range("I4").Value = 103.9136394821795
$wsArima.Range("B5").Value = 8.854566201129186
$wsArima.Range("C5").Value = 0.6161935655094513
$wsArima.Range("D5").Value = 0.6012842876856677
$wsArima.Range("E5").Value = 0.09242485518368102
$wsArima.Range("F5").Value = 2.767655747533855
$wsArima.Range("G5").Value = 5.282688612550318
$wsArima.Range("H5").Value = 103.6402991984743
$wsArima.Range("I5").Value = 10.1938039750713
$wsArima.Range("B6").Value = 0.03611511253830557
$wsArima.Range("C6").Value = 0.005987948643043344
$wsArima.Range("D6").Value = 0.005327286492916386
$wsArima.Range("E6").Value = 0.01469217224402913
$wsArima.Range("F6").Value = 0.005472580903409623
$wsArima.Range("G6").Value = 0.006323547321265265
$wsArima.Range("H6").Value = 0.02515096099664834
$wsArima.Range("I6").Value = 0.005214907209287058
$wsArima.Range("B7").Value = -6.315005592893259
$wsArima.Range("C7").Value = -0.06613684666586563
$wsArima.Range("D7").Value = -0.0397700188065121
$wsArima.Range("E7").Value = 0.002788150497838946
$wsArima.Range("F7").Value = -0.2585360032298287
$wsArima.Range("G7").Value = 0.2576701283004221
$wsArima.Range("H7").Value = -50.62132899077168
$wsArima.Range("I7").Value = -1.573897774164294
$wsArima.Range("B8").Value = 0.3556942499929322
$wsArima.Range("C8").Value = 0.9881939608847117
$wsArima.Range("D8").Value = 0.9840669136952048
$wsArima.Range("E8").Value = 0.994923496736837
$wsArima.Range("F8").Value = 0.9833319965041305
$wsArima.Range("G8").Value = 0.9960396910661393
$wsArima.Range("H8").Value = 0.4669624747535596
$wsArima.Range("I8").Value = 0.9932064909565916

# --- Sheet 2 ("lstm", formerly "rf") updated metric values ---
$wsRf.Range("B2").Value = 0.5447154471544715
$wsRf.Range("C2").Value = 0.5203252032520326
$wsRf.Range("D2").Value = 0.5338753387533876
$wsRf.Range("E2").Value = 0.5447154471544715
$wsRf.Range("F2").Value = 0.5203252032520326
$wsRf.Range("G2").Value = 0.5555555555555556
$wsRf.Range("H2").Value = 0.5094850948509485
$wsRf.Range("I2").Value = 0.5420054200542005
$wsRf.Range("B3").Value = 0.5352941176470588
$wsRf.Range("C3").Value = 0.5156695156695157
$wsRf.Range("D3").Value = 0.6875
$wsRf.Range("E3").Value = 0.5384615384615384
$wsRf.Range("F3").Value = 0
$wsRf.Range("G3").Value = 0.5230769230769231
$wsRf.Range("H3").Value = 0.5094850948509485
$wsRf.Range("I3").Value = 0.5909090909090909
$wsRf.Range("B4").Value = 5.697480619886242
$wsRf.Range("C4").Value = 3.971498249945091
$wsRf.Range("D4").Value = 3.00173506322719
$wsRf.Range("E4").Value = 0.03554718831105954
$wsRf.Range("F4").Value = 37.92529183524886
$wsRf.Range("G4").Value = 280.2300867827402
$wsRf.Range("H4").Value = 1207.817713495811
$wsRf.Range("I4").Value = 417.9626792648701
$wsRf.Range("B5").Value = 2.386939592843992
$wsRf.Range("C5").Value = 1.992861824097469
$wsRf.Range("D5").Value = 1.732551604780415
$wsRf.Range("E5").Value = 0.1885396200034877
$wsRf.Range("F5").Value = 6.158351389393825
$wsRf.Range("G5").Value = 16.7400742765001
$wsRf.Range("H5").Value = 34.75367194262803
$wsRf.Range("I5").Value = 20.44413557147551
$wsRf.Range("B6").Value = 0.01000047267455709
$wsRf.Range("C6").Value = 0.01928416841780949
$wsRf.Range("D6").Value = 0.01619349587629995
$wsRf.Range("E6").Value = 0.04942039684199519
$wsRf.Range("F6").Value = 0.01211904269064898
$wsRf.Range("G6").Value = 0.02031541258809129
$wsRf.Range("H6").Value = 0.0132353061016147
$wsRf.Range("I6").Value = 0.01120087218518798
$wsRf.Range("B7").Value = 0.840642003832148
$wsRf.Range("C7").Value = -0.1717958370715301
$wsRf.Range("D7").Value = 0.08035640147966379
$wsRf.Range("E7").Value = 0.03442302019744709
$wsRf.Range("F7").Value = 0.2496377092066809
$wsRf.Range("G7").Value = 1.987327513655996
$wsRf.Range("H7").Value = 4.541863897887025
$wsRf.Range("I7").Value = 3.662490803374831
$wsRf.Range("B8").Value = 0.9399396280246173
$wsRf.Range("C8").Value = 0.8507132023758059
$wsRf.Range("D8").Value = 0.8659494192009241
$wsRf.Range("E8").Value = 0.9549018176175083
$wsRf.Range("F8").Value = 0.9101012164328519
$wsRf.Range("G8").Value = 0.9426621649701536
$wsRf.Range("H8").Value = 0.9233740577815442
$wsRf.Range("I8").Value = 0.9742583587575426

